# Scheduled data refresh: update market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# for the rows whose crafted items had new Universalis price data.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 2397.7673
$ws.Range("I137").Value = 1951.2
$ws.Range("J137").Value = 3428.3076
$ws.Range("K137").Value = 5853.6
$ws.Range("L137").Value = 10284.9228
$ws.Range("M137").Value = -3303.6
$ws.Range("N137").Value = -15384.9228
# Row 138
$ws.Range("H138").Value = 4099.7104
$ws.Range("I138").Value = 3771.1304
$ws.Range("J138").Value = 4603.533
$ws.Range("K138").Value = 11313.3912
$ws.Range("L138").Value = 13810.599
$ws.Range("M138").Value = -6173.3912
$ws.Range("N138").Value = -24090.599

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6253366.5
$ws.Range("I45").Value = 14288491
$ws.Range("K45").Value = 14288491
$ws.Range("M45").Value = -14288114
# Row 61
$ws.Range("H61").Value = 8014.6333
$ws.Range("I61").Value = 7634.5264
$ws.Range("K61").Value = 7634.5264
$ws.Range("M61").Value = -7422.5264
# Row 74
$ws.Range("H74").Value = 6809.294
$ws.Range("I74").Value = 3518.1538
$ws.Range("K74").Value = 3518.1538
$ws.Range("M74").Value = -2644.1538
# Row 77
$ws.Range("H77").Value = 6809.294
$ws.Range("I77").Value = 3518.1538
$ws.Range("K77").Value = 17590.769
$ws.Range("M77").Value = -13222.769
# Row 132
$ws.Range("H132").Value = 2349.5715
$ws.Range("I132").Value = 1899.0938
$ws.Range("K132").Value = 5697.2814
$ws.Range("M132").Value = -3167.2814
# Row 136
$ws.Range("H136").Value = 8014.6333
$ws.Range("I136").Value = 7634.5264
$ws.Range("K136").Value = 22903.5792
$ws.Range("M136").Value = -20353.5792

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 755.75
$ws.Range("I22").Value = 569.8125
$ws.Range("J22").Value = 1499.5
$ws.Range("K22").Value = 569.8125
$ws.Range("L22").Value = 1499.5
$ws.Range("M22").Value = -396.8125
$ws.Range("N22").Value = -1845.5
# Row 94
$ws.Range("H94").Value = 1548.409
$ws.Range("J94").Value = 3206.2856
$ws.Range("L94").Value = 3206.2856
$ws.Range("N94").Value = -4108.2856
# Row 134
$ws.Range("H134").Value = 4582.3193
$ws.Range("I134").Value = 3316.2903
$ws.Range("J134").Value = 7035.25
$ws.Range("K134").Value = 9948.8709
$ws.Range("L134").Value = 21105.75
$ws.Range("M134").Value = -7413.8709
$ws.Range("N134").Value = -26175.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 135000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 135000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 135000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -135588
# Row 96
$ws.Range("H96").Value = 32159.8
$ws.Range("J96").Value = 32159.8
$ws.Range("L96").Value = 32159.8
$ws.Range("N96").Value = -37651.8
# Row 105
$ws.Range("H105").Value = 1834.3077
$ws.Range("I105").Value = 1834.3077
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1834.3077
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -87.30770000000007
$ws.Range("N105").ClearContents()
# Row 132
$ws.Range("H132").Value = 8401.666999999999
$ws.Range("I132").Value = 6860.0557
$ws.Range("J132").Value = 10714.083
$ws.Range("K132").Value = 20580.1671
$ws.Range("L132").Value = 32142.249
$ws.Range("M132").Value = -18050.1671
$ws.Range("N132").Value = -37202.249

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 2591
$ws.Range("I80").Value = 1981.6666
$ws.Range("J80").Value = 2852.1428
$ws.Range("K80").Value = 5944.9998
$ws.Range("L80").Value = 8556.428400000001
$ws.Range("M80").Value = -5008.9998
$ws.Range("N80").Value = -10428.4284
# Row 83
$ws.Range("H83").Value = 2591
$ws.Range("I83").Value = 1981.6666
$ws.Range("J83").Value = 2852.1428
$ws.Range("K83").Value = 17834.9994
$ws.Range("L83").Value = 25669.2852
$ws.Range("M83").Value = -13154.9994
$ws.Range("N83").Value = -35029.2852
# Row 98
$ws.Range("H98").Value = 246.81818
$ws.Range("I98").Value = 310.5
$ws.Range("J98").Value = 210.42857
$ws.Range("K98").Value = 931.5
$ws.Range("L98").Value = 631.28571
$ws.Range("M98").Value = 566.5
$ws.Range("N98").Value = -3627.28571
# Row 122
$ws.Range("H122").Value = 866.2857
$ws.Range("I122").Value = 485
$ws.Range("K122").Value = 4365
$ws.Range("M122").Value = -1915

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 132449
$ws.Range("I35").Value = 39000
$ws.Range("J35").Value = 225898
$ws.Range("K35").Value = 39000
$ws.Range("L35").Value = 225898
$ws.Range("M35").Value = -38702
$ws.Range("N35").Value = -226494
# Row 39
$ws.Range("H39").Value = 117333.336
$ws.Range("J39").Value = 146000
$ws.Range("L39").Value = 146000
$ws.Range("N39").Value = -147064
# Row 70
$ws.Range("H70").Value = 4991.6
$ws.Range("I70").Value = 4443.5
$ws.Range("K70").Value = 4443.5
$ws.Range("M70").Value = -4173.5
# Row 73
$ws.Range("H73").Value = 4991.6
$ws.Range("I73").Value = 4443.5
$ws.Range("K73").Value = 4443.5
$ws.Range("M73").Value = -3507.5
# Row 132
$ws.Range("H132").Value = 4926.306
$ws.Range("I132").Value = 4090.15
$ws.Range("K132").Value = 12270.45
$ws.Range("M132").Value = -9740.450000000001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3532.606
$ws.Range("I122").Value = 2004.2609
$ws.Range("K122").Value = 6012.7827
$ws.Range("M122").Value = -3562.7827
# Row 127
$ws.Range("H127").Value = 93000
$ws.Range("J127").Value = 93000
$ws.Range("L127").Value = 93000
$ws.Range("N127").Value = -102920
# Row 132
$ws.Range("H132").Value = 6663.4473
$ws.Range("I132").Value = 5985.852
$ws.Range("J132").Value = 8326.637000000001
$ws.Range("K132").Value = 17957.556
$ws.Range("L132").Value = 24979.911
$ws.Range("M132").Value = -15427.556
$ws.Range("N132").Value = -30039.911
# Row 136
$ws.Range("H136").Value = 8528.799999999999
$ws.Range("I136").Value = 3949.111
$ws.Range("K136").Value = 11847.333
$ws.Range("M136").Value = -9297.332999999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3404.0833
$ws.Range("I126").Value = 2085
$ws.Range("J126").Value = 9999.5
$ws.Range("K126").Value = 6255
$ws.Range("L126").Value = 29998.5
$ws.Range("M126").Value = -3785
$ws.Range("N126").Value = -34938.5
# Row 132
$ws.Range("H132").Value = 2662.2188
$ws.Range("I132").Value = 1989.7
$ws.Range("J132").Value = 12750
$ws.Range("K132").Value = 5969.1
$ws.Range("L132").Value = 38250
$ws.Range("M132").Value = -3439.1
$ws.Range("N132").Value = -43310
# Row 136
$ws.Range("H136").Value = 7023.533
$ws.Range("I136").Value = 5929.222
$ws.Range("K136").Value = 17787.666
$ws.Range("M136").Value = -15237.666
